$wb = $excel.ActiveWorkbook

# --- Side3 sheet (sheet4.xml) ---
$wsSide3 = $wb.Worksheets.Item("Side3")
$wsSide3.Activate()
# Clear the swap-in character skill rows (14-17), columns B:C, and remove the stray D14 value
$wsSide3.Range("B14:C17").ClearContents()
$wsSide3.Range("D14").ClearContents()
# update the active selection to C23
$wsSide3.Range("C23").Select()

# --- Settings sheet (sheet1.xml) ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Activate()
# B11: "master clothes is swap-in costume" flag 1 -> 0 (bug fix for character checking)
$wsSettings.Range("B11").Value = 0
# update the active selection to B7
$wsSettings.Range("B7").Select()
